# Update lexical diversity values to count as a proportion of all tokens,
# not just word tokens (per-document values across the 2009-2018 sheets, and
# the recomputed count/mean/std/min/25%/50%/75%/max stats on the Summary sheet).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2009")
$ws.Range("B2").Value = 0.3512195121951219
$ws.Range("B3").Value = 0.5133079847908745
$ws.Range("B4").Value = 0.4244444444444445
$ws.Range("B5").Value = 0.431980906921241
$ws.Range("B6").Value = 0.3670309653916211
$ws.Range("B7").Value = 0.4499089253187614
$ws.Range("B8").Value = 0.3818722139673105
$ws.Range("B9").Value = 0.4292763157894737
$ws.Range("B10").Value = 0.4164358264081256
$ws.Range("B11").Value = 0.456140350877193
$ws.Range("B12").Value = 0.3336820083682008

$ws = $wb.Worksheets.Item("2018")
$ws.Range("B2").Value = 0.4726890756302521
$ws.Range("B3").Value = 0.5177111716621253
$ws.Range("B4").Value = 0.2869222096956032
$ws.Range("B5").Value = 0.4968553459119497
$ws.Range("B6").Value = 0.4626647144948756
$ws.Range("B7").Value = 0.3132036847492323
$ws.Range("B8").Value = 0.3814634146341463
$ws.Range("B9").Value = 0.4175182481751825
$ws.Range("B10").Value = 0.4386459802538787
$ws.Range("B11").Value = 0.426497277676951
$ws.Range("B12").Value = 0.5241379310344828
$ws.Range("B13").Value = 0.4970414201183432
$ws.Range("B14").Value = 0.3629032258064516
$ws.Range("B15").Value = 0.5244299674267101
$ws.Range("B16").Value = 0.5586734693877551
$ws.Range("B17").Value = 0.4417344173441735
$ws.Range("B18").Value = 0.3772582359192349
$ws.Range("B19").Value = 0.7701149425287356
$ws.Range("B20").Value = 0.4686411149825784
$ws.Range("B21").Value = 0.5910364145658263
$ws.Range("B22").Value = 0.6217391304347826
$ws.Range("B23").Value = 0.6242038216560509
$ws.Range("B24").Value = 0.4494845360824742
$ws.Range("B25").Value = 0.3517179023508137
$ws.Range("B26").Value = 0.5119760479041916
$ws.Range("B27").Value = 0.6568047337278107
$ws.Range("B28").Value = 0.4415156507413509
$ws.Range("B29").Value = 0.4568393094289508
$ws.Range("B30").Value = 0.640625
$ws.Range("B31").Value = 0.3630252100840336
$ws.Range("B32").Value = 0.4320137693631669
$ws.Range("B33").Value = 0.3186440677966101
$ws.Range("B34").Value = 0.4771784232365145
$ws.Range("B35").Value = 0.5719063545150501
$ws.Range("B36").Value = 0.3818965517241379
$ws.Range("B37").Value = 0.3678025851938895
$ws.Range("B38").Value = 0.4985422740524781
$ws.Range("B39").Value = 0.4693140794223827
$ws.Range("B40").Value = 0.4674922600619195
$ws.Range("B41").Value = 0.6118721461187214
$ws.Range("B42").Value = 0.4217391304347826
$ws.Range("B43").Value = 0.4053763440860215
$ws.Range("B44").Value = 0.5077399380804953
$ws.Range("B45").Value = 0.489247311827957
$ws.Range("B46").Value = 0.3582718651211802
$ws.Range("B47").Value = 0.3877840909090909
$ws.Range("B48").Value = 0.4120879120879121
$ws.Range("B49").Value = 0.3836363636363636
$ws.Range("B50").Value = 0.4593175853018373
$ws.Range("B51").Value = 0.467032967032967
$ws.Range("B52").Value = 0.4423076923076923
$ws.Range("B53").Value = 0.4902912621359223
$ws.Range("B54").Value = 0.6161616161616161
$ws.Range("B55").Value = 0.3417493237150586
$ws.Range("B56").Value = 0.3987951807228916
$ws.Range("B57").Value = 0.4666666666666667
$ws.Range("B58").Value = 0.4508670520231214
$ws.Range("B59").Value = 0.5641025641025641
$ws.Range("B60").Value = 0.3903903903903904
$ws.Range("B61").Value = 0.4655581947743468
$ws.Range("B62").Value = 0.464968152866242
$ws.Range("B63").Value = 0.4695121951219512
$ws.Range("B64").Value = 0.6422764227642277
$ws.Range("B65").Value = 0.6141078838174274
$ws.Range("B66").Value = 0.5916666666666667
$ws.Range("B67").Value = 0.5508474576271186
$ws.Range("B68").Value = 0.4616541353383459
$ws.Range("B69").Value = 0.5242165242165242
$ws.Range("B70").Value = 0.6643356643356644
$ws.Range("B71").Value = 0.4462151394422311
$ws.Range("B72").Value = 0.4328899637243047
$ws.Range("B73").Value = 0.6357615894039735
$ws.Range("B74").Value = 0.4918414918414918
$ws.Range("B75").Value = 0.5547945205479452
$ws.Range("B76").Value = 0.4163179916317992
$ws.Range("B77").Value = 0.6043613707165109
$ws.Range("B78").Value = 0.433048433048433
$ws.Range("B79").Value = 0.3616298811544992
$ws.Range("B80").Value = 0.3283292978208233
$ws.Range("B81").Value = 0.4672131147540984
$ws.Range("B82").Value = 0.4730831973898858
$ws.Range("B83").Value = 0.5080645161290323
$ws.Range("B84").Value = 0.4763313609467456
$ws.Range("B85").Value = 0.5543859649122806
$ws.Range("B86").Value = 0.4501557632398754

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("C2").Value = 0.4141181322247607
$ws.Range("D2").Value = 0.0521819271486032
$ws.Range("E2").Value = 0.3336820083682008
$ws.Range("F2").Value = 0.3744515896794658
$ws.Range("G2").Value = 0.4244444444444445
$ws.Range("H2").Value = 0.4409449161200012
$ws.Range("I2").Value = 0.5133079847908745
$ws.Range("C3").Value = 0.4391663905643776
$ws.Range("D3").Value = 0.07233673996480966
$ws.Range("E3").Value = 0.2951200619674671
$ws.Range("F3").Value = 0.3907113960146986
$ws.Range("G3").Value = 0.4335922330097087
$ws.Range("H3").Value = 0.4737946998722861
$ws.Range("I3").Value = 0.547945205479452
$ws.Range("C4").Value = 0.4496205507584851
$ws.Range("D4").Value = 0.0858020596468571
$ws.Range("E4").Value = 0.3452544704264099
$ws.Range("F4").Value = 0.3877745304303952
$ws.Range("G4").Value = 0.4351732991014121
$ws.Range("H4").Value = 0.4781528715512211
$ws.Range("I4").Value = 0.7583333333333333
$ws.Range("C5").Value = 0.4316655795449358
$ws.Range("D5").Value = 0.08173450336538521
$ws.Range("E5").Value = 0.2432018383761011
$ws.Range("F5").Value = 0.3987924725022323
$ws.Range("G5").Value = 0.4157274875597994
$ws.Range("H5").Value = 0.4843366093366094
$ws.Range("I5").Value = 0.6017699115044248
$ws.Range("C6").Value = 0.4177431768965846
$ws.Range("D6").Value = 0.09582027019874044
$ws.Range("E6").Value = 0.2681692732290709
$ws.Range("F6").Value = 0.3471676843241325
$ws.Range("G6").Value = 0.418265628819962
$ws.Range("H6").Value = 0.4672378979525648
$ws.Range("I6").Value = 0.6753246753246753
$ws.Range("C7").Value = 0.4584010266073282
$ws.Range("D7").Value = 0.100004732064552
$ws.Range("E7").Value = 0.2825443786982249
$ws.Range("F7").Value = 0.4141025641025641
$ws.Range("G7").Value = 0.4393063583815029
$ws.Range("H7").Value = 0.48
$ws.Range("I7").Value = 0.7045454545454546
$ws.Range("C8").Value = 0.4239425005351002
$ws.Range("D8").Value = 0.08101709676570246
$ws.Range("E8").Value = 0.220218319886094
$ws.Range("F8").Value = 0.383687562856185
$ws.Range("G8").Value = 0.4421875
$ws.Range("H8").Value = 0.4761144113651539
$ws.Range("I8").Value = 0.5397727272727273
$ws.Range("C9").Value = 0.4390697247199625
$ws.Range("D9").Value = 0.08800602947157546
$ws.Range("E9").Value = 0.2385786802030457
$ws.Range("F9").Value = 0.385987460815047
$ws.Range("G9").Value = 0.4380403458213257
$ws.Range("H9").Value = 0.4961137798888109
$ws.Range("I9").Value = 0.6575342465753424
$ws.Range("C10").Value = 0.497163422256954
$ws.Range("D10").Value = 0.09470659441718095
$ws.Range("E10").Value = 0.3114134542705971
$ws.Range("F10").Value = 0.432684157805517
$ws.Range("G10").Value = 0.4943181818181818
$ws.Range("H10").Value = 0.5440259971509971
$ws.Range("I10").Value = 0.9666666666666667
$ws.Range("C11").Value = 0.4774337443161269
$ws.Range("D11").Value = 0.09349725745545946
$ws.Range("E11").Value = 0.2869222096956032
$ws.Range("F11").Value = 0.4175182481751825
$ws.Range("G11").Value = 0.467032967032967
$ws.Range("H11").Value = 0.5242165242165242
$ws.Range("I11").Value = 0.7701149425287356

$ws = $wb.Worksheets.Item("2010")
$ws.Range("B2").Value = 0.3819691577698695
$ws.Range("B3").Value = 0.4271844660194175
$ws.Range("B4").Value = 0.44
$ws.Range("B5").Value = 0.5321100917431193
$ws.Range("B6").Value = 0.5437352245862884
$ws.Range("B7").Value = 0.372310570626754
$ws.Range("B8").Value = 0.547945205479452
$ws.Range("B9").Value = 0.3713620488940629
$ws.Range("B10").Value = 0.4189349112426036
$ws.Range("B11").Value = 0.4169381107491857
$ws.Range("B12").Value = 0.4572649572649573
$ws.Range("B13").Value = 0.4675925925925926
$ws.Range("B14").Value = 0.4758620689655172
$ws.Range("B15").Value = 0.2951200619674671

$ws = $wb.Worksheets.Item("2011")
$ws.Range("B2").Value = 0.3452544704264099
$ws.Range("B3").Value = 0.5643564356435643
$ws.Range("B4").Value = 0.3742116327960757
$ws.Range("B5").Value = 0.5497630331753555
$ws.Range("B6").Value = 0.3744725738396624
$ws.Range("B7").Value = 0.3830472103004292
$ws.Range("B8").Value = 0.4692737430167598
$ws.Range("B9").Value = 0.3892215568862276
$ws.Range("B10").Value = 0.3965267727930535
$ws.Range("B11").Value = 0.364
$ws.Range("B12").Value = 0.4683098591549296
$ws.Range("B13").Value = 0.5012594458438288
$ws.Range("B14").Value = 0.3863275039745628
$ws.Range("B15").Value = 0.4391304347826087
$ws.Range("B16").Value = 0.3511214230471771
$ws.Range("B17").Value = 0.4351732991014121
$ws.Range("B18").Value = 0.4334140435835351
$ws.Range("B19").Value = 0.7583333333333333
$ws.Range("B20").Value = 0.4813953488372093
$ws.Range("B21").Value = 0.4992101105845181
$ws.Range("B22").Value = 0.4112792297111417
$ws.Range("B23").Value = 0.4421364985163205
$ws.Range("B24").Value = 0.4299424184261036
$ws.Range("B25").Value = 0.474910394265233
$ws.Range("B26").Value = 0.5373134328358209
$ws.Range("B27").Value = 0.4192307692307692
$ws.Range("B28").Value = 0.461139896373057

$ws = $wb.Worksheets.Item("2012")
$ws.Range("B2").Value = 0.4276527331189711
$ws.Range("B3").Value = 0.389428263214671
$ws.Range("B4").Value = 0.4044117647058824
$ws.Range("B5").Value = 0.5
$ws.Range("B6").Value = 0.539568345323741
$ws.Range("B7").Value = 0.6017699115044248
$ws.Range("B8").Value = 0.4081632653061225
$ws.Range("B9").Value = 0.502283105022831
$ws.Range("B10").Value = 0.2432018383761011
$ws.Range("B11").Value = 0.3882352941176471
$ws.Range("B12").Value = 0.4219178082191781
$ws.Range("B13").Value = 0.3558504221954162
$ws.Range("B14").Value = 0.4791154791154791
$ws.Range("B15").Value = 0.4095371669004207
$ws.Range("B16").Value = 0.4336
$ws.Range("B17").Value = 0.4019138755980861

$ws = $wb.Worksheets.Item("2013")
$ws.Range("B2").Value = 0.2681692732290709
$ws.Range("B3").Value = 0.2681692732290709
$ws.Range("B4").Value = 0.4179620034542315
$ws.Range("B5").Value = 0.3748753738783649
$ws.Range("B6").Value = 0.6753246753246753
$ws.Range("B7").Value = 0.3400383141762452
$ws.Range("B8").Value = 0.3452593917710197
$ws.Range("B9").Value = 0.396640826873385
$ws.Range("B10").Value = 0.4545454545454545
$ws.Range("B11").Value = 0.3528925619834711
$ws.Range("B12").Value = 0.446078431372549
$ws.Range("B13").Value = 0.4691943127962085
$ws.Range("B14").Value = 0.5109649122807017
$ws.Range("B15").Value = 0.5508982035928144
$ws.Range("B16").Value = 0.34
$ws.Range("B17").Value = 0.358053302433372
$ws.Range("B18").Value = 0.3619489559164733
$ws.Range("B19").Value = 0.4011049723756906
$ws.Range("B20").Value = 0.4702602230483272
$ws.Range("B21").Value = 0.5388127853881278
$ws.Range("B22").Value = 0.479108635097493
$ws.Range("B23").Value = 0.4371002132196162
$ws.Range("B24").Value = 0.3000931966449208
$ws.Range("B25").Value = 0.3000931966449208
$ws.Range("B26").Value = 0.3000931966449208
$ws.Range("B27").Value = 0.354066985645933
$ws.Range("B28").Value = 0.6358695652173914
$ws.Range("B29").Value = 0.4569983136593592
$ws.Range("B30").Value = 0.4185692541856925
$ws.Range("B31").Value = 0.3392996108949417
$ws.Range("B32").Value = 0.449645390070922
$ws.Range("B33").Value = 0.4571428571428571
$ws.Range("B34").Value = 0.4726256983240223
$ws.Range("B35").Value = 0.4613686534216336

$ws = $wb.Worksheets.Item("2014")
$ws.Range("B2").Value = 0.4490644490644491
$ws.Range("B3").Value = 0.7045454545454546
$ws.Range("B4").Value = 0.4183535762483131
$ws.Range("B5").Value = 0.4393063583815029
$ws.Range("B6").Value = 0.6578947368421053
$ws.Range("B7").Value = 0.4141025641025641
$ws.Range("B8").Value = 0.502283105022831
$ws.Range("B9").Value = 0.4582043343653251
$ws.Range("B10").Value = 0.4344827586206896
$ws.Range("B11").Value = 0.3755615453728661
$ws.Range("B12").Value = 0.2825443786982249
$ws.Range("B13").Value = 0.4102564102564102
$ws.Range("B14").Value = 0.4274891774891775
$ws.Range("B15").Value = 0.4739776951672863
$ws.Range("B16").Value = 0.48
$ws.Range("B17").Value = 0.3617021276595745
$ws.Range("B18").Value = 0.5030487804878049

$ws = $wb.Worksheets.Item("2015")
$ws.Range("B2").Value = 0.4467005076142132
$ws.Range("B3").Value = 0.5397727272727273
$ws.Range("B4").Value = 0.4704301075268817
$ws.Range("B5").Value = 0.4850498338870432
$ws.Range("B6").Value = 0.4170692431561997
$ws.Range("B7").Value = 0.2853817504655494
$ws.Range("B8").Value = 0.4438305709023941
$ws.Range("B9").Value = 0.5128205128205128
$ws.Range("B10").Value = 0.4365620736698499
$ws.Range("B11").Value = 0.5170454545454546
$ws.Range("B12").Value = 0.3232791956689868
$ws.Range("B13").Value = 0.4693446088794926
$ws.Range("B14").Value = 0.321954484605087
$ws.Range("B15").Value = 0.391578947368421
$ws.Range("B16").Value = 0.5145985401459854
$ws.Range("B17").Value = 0.4421875
$ws.Range("B18").Value = 0.4525139664804469
$ws.Range("B19").Value = 0.3757961783439491
$ws.Range("B20").Value = 0.4265927977839335
$ws.Range("B21").Value = 0.220218319886094
$ws.Range("B22").Value = 0.440251572327044
$ws.Range("B23").Value = 0.3358999037536092
$ws.Range("B24").Value = 0.4817987152034261

$ws = $wb.Worksheets.Item("2016")
$ws.Range("B2").Value = 0.3315696649029982
$ws.Range("B3").Value = 0.4380403458213257
$ws.Range("B4").Value = 0.3640040444893832
$ws.Range("B5").Value = 0.5915492957746479
$ws.Range("B6").Value = 0.5646551724137931
$ws.Range("B7").Value = 0.4741784037558686
$ws.Range("B8").Value = 0.2971788551081895
$ws.Range("B9").Value = 0.3758911211924822
$ws.Range("B10").Value = 0.4602803738317757
$ws.Range("B11").Value = 0.5328947368421053
$ws.Range("B12").Value = 0.4305835010060362
$ws.Range("B13").Value = 0.3884555382215288
$ws.Range("B14").Value = 0.3952225841476656
$ws.Range("B15").Value = 0.4059787849566056
$ws.Range("B16").Value = 0.3748427672955975
$ws.Range("B17").Value = 0.4261168384879725
$ws.Range("B18").Value = 0.4939516129032258
$ws.Range("B19").Value = 0.3944827586206897
$ws.Range("B20").Value = 0.4608150470219436
$ws.Range("B21").Value = 0.4395424836601307
$ws.Range("B22").Value = 0.5067873303167421
$ws.Range("B23").Value = 0.4236526946107784
$ws.Range("B24").Value = 0.3107658157602664
$ws.Range("B25").Value = 0.4308093994778068
$ws.Range("B26").Value = 0.5228070175438596
$ws.Range("B27").Value = 0.2385786802030457
$ws.Range("B28").Value = 0.4897959183673469
$ws.Range("B29").Value = 0.3848275862068966
$ws.Range("B30").Value = 0.2931323283082077
$ws.Range("B31").Value = 0.4948453608247423
$ws.Range("B32").Value = 0.5799086757990868
$ws.Range("B33").Value = 0.4973821989528796
$ws.Range("B34").Value = 0.5433526011560693
$ws.Range("B35").Value = 0.3871473354231975
$ws.Range("B36").Value = 0.4656488549618321
$ws.Range("B37").Value = 0.4062038404726735
$ws.Range("B38").Value = 0.4593301435406699
$ws.Range("B39").Value = 0.4438502673796791
$ws.Range("B40").Value = 0.6575342465753424
$ws.Range("B41").Value = 0.5018587360594795
$ws.Range("B42").Value = 0.4639498432601881
$ws.Range("B43").Value = 0.3808080808080808
$ws.Range("B44").Value = 0.3162692847124824
$ws.Range("B45").Value = 0.2942260442260442
$ws.Range("B46").Value = 0.4237588652482269
$ws.Range("B47").Value = 0.5536723163841808
$ws.Range("B48").Value = 0.5251396648044693

$ws = $wb.Worksheets.Item("2017")
$ws.Range("B2").Value = 0.5012106537530266
$ws.Range("B3").Value = 0.4512922465208747
$ws.Range("B4").Value = 0.4623287671232877
$ws.Range("B5").Value = 0.4125098970704671
$ws.Range("B6").Value = 0.400208986415883
$ws.Range("B7").Value = 0.3114134542705971
$ws.Range("B8").Value = 0.4643478260869565
$ws.Range("B9").Value = 0.5642023346303502
$ws.Range("B10").Value = 0.4368421052631579
$ws.Range("B11").Value = 0.5329949238578681
$ws.Range("B12").Value = 0.4730538922155689
$ws.Range("B13").Value = 0.3735035913806863
$ws.Range("B14").Value = 0.5254237288135594
$ws.Range("B15").Value = 0.4175392670157068
$ws.Range("B16").Value = 0.4328097731239092
$ws.Range("B17").Value = 0.4237588652482269
$ws.Range("B18").Value = 0.5195530726256983
$ws.Range("B19").Value = 0.4327323162274618
$ws.Range("B20").Value = 0.4780793319415449
$ws.Range("B21").Value = 0.4961636828644501
$ws.Range("B22").Value = 0.5168539325842697
$ws.Range("B23").Value = 0.3814589665653496
$ws.Range("B24").Value = 0.5509433962264151
$ws.Range("B25").Value = 0.5140562248995983
$ws.Range("B26").Value = 0.3772455089820359
$ws.Range("B27").Value = 0.3876651982378855
$ws.Range("B28").Value = 0.624390243902439
$ws.Range("B29").Value = 0.4809619238476954
$ws.Range("B30").Value = 0.5882352941176471
$ws.Range("B31").Value = 0.4161220043572985
$ws.Range("B32").Value = 0.5201793721973094
$ws.Range("B33").Value = 0.4733893557422969
$ws.Range("B34").Value = 0.4140893470790378
$ws.Range("B35").Value = 0.6379310344827587
$ws.Range("B36").Value = 0.5462962962962963
$ws.Range("B37").Value = 0.5256410256410257
$ws.Range("B38").Value = 0.5318471337579618
$ws.Range("B39").Value = 0.5144927536231884
$ws.Range("B40").Value = 0.5542168674698795
$ws.Range("B41").Value = 0.504
$ws.Range("B42").Value = 0.527972027972028
$ws.Range("B43").Value = 0.3675675675675676
$ws.Range("B44").Value = 0.508
$ws.Range("B45").Value = 0.5943775100401606
$ws.Range("B46").Value = 0.5714285714285714
$ws.Range("B47").Value = 0.5346534653465347
$ws.Range("B48").Value = 0.4139344262295082
$ws.Range("B49").Value = 0.52734375
$ws.Range("B50").Value = 0.554140127388535
$ws.Range("B51").Value = 0.5
$ws.Range("B52").Value = 0.5464684014869888
$ws.Range("B53").Value = 0.5754385964912281
$ws.Range("B54").Value = 0.5225464190981433
$ws.Range("B55").Value = 0.6772486772486772
$ws.Range("B56").Value = 0.4135021097046414
$ws.Range("B57").Value = 0.4105461393596987
$ws.Range("B58").Value = 0.4863523573200992
$ws.Range("B59").Value = 0.4879725085910653
$ws.Range("B60").Value = 0.4618181818181818
$ws.Range("B61").Value = 0.4009779951100245
$ws.Range("B62").Value = 0.3818443804034582
$ws.Range("B63").Value = 0.9666666666666667
$ws.Range("B64").Value = 0.4362745098039216
$ws.Range("B65").Value = 0.3652058432934927
$ws.Range("B66").Value = 0.4774193548387097
$ws.Range("B67").Value = 0.4643237486687966
$ws.Range("B68").Value = 0.4609375
$ws.Range("B69").Value = 0.3888070692194404
$ws.Range("B70").Value = 0.4519230769230769
$ws.Range("B71").Value = 0.4228187919463087
$ws.Range("B72").Value = 0.4527027027027027
$ws.Range("B73").Value = 0.4968553459119497
$ws.Range("B74").Value = 0.4960212201591512
$ws.Range("B75").Value = 0.3568129330254042
$ws.Range("B76").Value = 0.531986531986532
$ws.Range("B77").Value = 0.6081081081081081
$ws.Range("B78").Value = 0.462962962962963
$ws.Range("B79").Value = 0.5294117647058824
$ws.Range("B80").Value = 0.6343283582089553
$ws.Range("B81").Value = 0.6086956521739131
$ws.Range("B82").Value = 0.5065963060686016
$ws.Range("B83").Value = 0.5286783042394015
$ws.Range("B84").Value = 0.4558303886925795
$ws.Range("B85").Value = 0.4158800666296502
$ws.Range("B86").Value = 0.5405405405405406
$ws.Range("B87").Value = 0.5333333333333333
$ws.Range("B88").Value = 0.4943181818181818
$ws.Range("B89").Value = 0.4052044609665427
$ws.Range("B90").Value = 0.4787878787878788
$ws.Range("B91").Value = 0.6578947368421053
$ws.Range("B92").Value = 0.4943181818181818
$ws.Range("B93").Value = 0.4321678321678322
$ws.Range("B94").Value = 0.3890865954922895
$ws.Range("B95").Value = 0.592814371257485
$ws.Range("B96").Value = 0.4389610389610389
$ws.Range("B97").Value = 0.4718309859154929
$ws.Range("B98").Value = 0.4598698481561822
$ws.Range("B99").Value = 0.5975103734439834
$ws.Range("B100").Value = 0.6139240506329114
$ws.Range("B101").Value = 0.5555555555555556
$ws.Range("B102").Value = 0.9
$ws.Range("B103").Value = 0.5740740740740741
$ws.Range("B104").Value = 0.6507936507936508
$ws.Range("B105").Value = 0.488
$ws.Range("B106").Value = 0.5551020408163265
$ws.Range("B107").Value = 0.572347266881029
$ws.Range("B108").Value = 0.6644295302013423
$ws.Range("B109").Value = 0.5625
$ws.Range("B110").Value = 0.4943181818181818
$ws.Range("B111").Value = 0.4818652849740933
$ws.Range("B112").Value = 0.4292035398230089
$ws.Range("B113").Value = 0.5366568914956011
$ws.Range("B114").Value = 0.3340060544904137
$ws.Range("B115").Value = 0.5432692307692307
$ws.Range("B116").Value = 0.5464285714285714
$ws.Range("B117").Value = 0.3210702341137124
$ws.Range("B118").Value = 0.5347721822541966
$ws.Range("B119").Value = 0.3400167084377611
$ws.Range("B120").Value = 0.5318518518518518
$ws.Range("B121").Value = 0.4822006472491909
$ws.Range("B122").Value = 0.5898617511520737
$ws.Range("B123").Value = 0.5521885521885522
$ws.Range("B124").Value = 0.5812807881773399
$ws.Range("B125").Value = 0.390194075587334
$ws.Range("B126").Value = 0.536
$ws.Range("B127").Value = 0.6694214876033058
$ws.Range("B128").Value = 0.4759358288770054
$ws.Range("B129").Value = 0.4735099337748344
$ws.Range("B130").Value = 0.4398340248962656
$ws.Range("B131").Value = 0.4325396825396826
$ws.Range("B132").Value = 0.4836879432624113
$ws.Range("B133").Value = 0.5032258064516129
$ws.Range("B134").Value = 0.5721153846153846
$ws.Range("B135").Value = 0.4159090909090909
$ws.Range("B136").Value = 0.4892857142857143
$ws.Range("B137").Value = 0.3501646542261251
$ws.Range("B138").Value = 0.3513513513513514
$ws.Range("B139").Value = 0.5357142857142857
$ws.Range("B140").Value = 0.5302593659942363
$ws.Range("B141").Value = 0.3780104712041885

